$wb = $excel.ActiveWorkbook

# --- 1. Status text change: "Ready for handoff" -> "In Translation" ---
# This status string is shared across the Overview sheet (columns E/F,
# one per target locale) and each locale sheet's "Status" column (C).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # Cast to [string] first: some cell values (e.g. the literal text
        # "True") otherwise get PowerShell's automatic boolean coercion,
        # which would make "True" -eq "Ready for handoff" spuriously true.
        if ([string]$cell.Value2 -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

# --- 2. Narrow the "Status" columns that held the long status text ---
# Overview sheet: zh-cn (E) and de-de (F) status columns.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# Per-locale sheets: Status column (C).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
